# Template (save as a csv).xlsx update
#  - Replace the "FileExtension" column header with "FileMatching"
#  - Rewrite the instructional/example rows in column D to describe the
#    new regex-based matching behaviour (instead of the old wildcard/-like
#    behaviour), and reorder a couple of rows
#  - Turn the last instructional row into a hyperlink pointing at a page
#    that explains/tests regex patterns
#  - Update the active-cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - "delete this row" notice (unchanged text, just re-levelled in the
# shared string table by the engine)
$ws.Range("A1").Value = "DELETE THIS ROW WHEN YOU SAVE:  The default process is backup, so when you run a backup, it will read from the ""Source"" column and put it into the ""Destination"" column."

# Row 2 - header row: FileExtension -> FileMatching
$ws.Range("A2").Value = "Source"
$ws.Range("B2").Value = "Destination"
$ws.Range("C2").Value = "Description"
$ws.Range("D2").Value = "FileMatching"

# Rows 3-11: new/updated explanatory text describing the regex matching
$ws.Range("D3").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Using a regex pattern, you can fetch the desired files seperated by a '/' for every pattern/entry you want to match to."
$ws.Range("D4").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  or you can leave it blank to get ALL the contents of the folder.  This is case insensitive so you don't need to worry about capital letters."
$ws.Range("D5").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Here are some examples all you have to do is to remove the '' around the text."
$ws.Range("D6").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: '.txt$' will get all files that end in '.txt' inside it's name.  The $ means the end of the string."
$ws.Range("D7").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: '^text' will get all files that start with 'text' inside it's name.  The ^ means at the start of the string."
$ws.Range("D8").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: 'text' will get all files that contains 'text' inside it's name."
$ws.Range("D9").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: '^text(?:\N)*.txt$' will get all files that start with 'text' with anything between the start and the end, even nothing, and ends with '.txt' "
$ws.Range("D10").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  As shown by the last example you can combine them to make complex patterns such as the following"
$ws.Range("D11").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: '^test (?:\N)*file(?:\N)* name\.txt$' will get all files that start with 'start ', and somewhere in the middle has 'file', and ends with  ' name.txt'"

# Row 12: now a hyperlink pointing at a regex reference/tester site
$ws.Range("D12").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  A good place to test and learn about regex patterns can be found here"
$ws.Hyperlinks.Add($ws.Range("D12"), "https://regex101.com/", "", "", "REMOVE\REPLACE ME WHEN YOU SAVE:  A good place to test and learn about regex patterns can be found here")

# Update the saved selection to D2 (matches the new active cell in the diff)
$ws.Range("D2").Select()
